# Generate Report for Handback
#
# This mirrors a re-run of the handback status report generation: the
# "Latest HO Xliff Generate Date" on the Overview sheet plus the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" columns
# on the per-language sheets are refreshed for the 0f64f7fa... file with
# newer timestamps produced by the latest handback pass.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 02:57:42"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 02:57:37"
$wsZhCn.Range("K2").Value = "2016-09-01 02:57:54"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 02:57:42"
$wsDeDe.Range("K2").Value = "2016-09-01 02:58:06"
